# Duble magyar módszer debuggolva, üres táblázaton jól működik.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header counters -------------------------------------------------
# "Napok száma" (number of days) drops from 10 to 7
$ws.Range("B1").Value = 7

# New helper totals next to the day / night shift headers
$ws.Range("C1").Formula = "=SUM(B4:B10)"
$ws.Range("C2").Formula = "=SUM(C4:C10)"

# --- Per-person day / night shift counts (B4:C10) ---------------------
$dayCounts = @{ 4 = 2; 5 = 2; 6 = 2; 7 = 1; 8 = 3; 9 = 2; 10 = 2 }
$nightCounts = @{ 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1 }

foreach ($row in 4..10) {
    $ws.Cells.Item($row, 2).Value = $dayCounts[$row]
    $ws.Cells.Item($row, 3).Value = $nightCounts[$row]
}

# --- Clear the old "n" / "é" / "N" / "É" shift-mark cells -------------
$clearCells = @(
    "D4","E4","F4","H4","M4",
    "D5","M5",
    "D6","E6","J6","K6","M6",
    "E7","H7","I7","M7",
    "M8",
    "D9","M9",
    "D10","F10","J10","K10"
)
foreach ($c in $clearCells) {
    $ws.Range($c).Value = ""
}

# --- Selection moves from E4 to B2 ------------------------------------
$ws.Range("B2").Select()
